$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths to match the new data.
# Excel's ColumnWidth COM property pads the stored OOXML <col width> by
# ~5px (0.8333... chars at the default font), so subtract that padding
# here to land on the exact target "width" values after round-tripping.
$padding = 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 55 - $padding
$ws.Columns.Item(3).ColumnWidth = 27 - $padding
$ws.Columns.Item(4).ColumnWidth = 28 - $padding
$ws.Columns.Item(5).ColumnWidth = 10 - $padding
$ws.Columns.Item(6).ColumnWidth = 16 - $padding
$ws.Columns.Item(7).ColumnWidth = 16 - $padding
$ws.Columns.Item(8).ColumnWidth = 38 - $padding

# New scraped opportunity rows to append (rows 2-6)
$data = @(
    @("1330737", "https://aiesec.org/opportunity/global-talent/1330737", "Web Developer Intern", "Phagwara, Punjab, India", "No", "2 applicants", "3 - 6 Months", "GNA University"),
    @("1330640", "https://aiesec.org/opportunity/global-talent/1330640", "Guest Relations Officer", "Galle, Sri Lanka", "No", "3 applicants", "3 - 6 Months", "Sino Lanka Hotels Colombo (Pvt) Ltd"),
    @("1329895", "https://aiesec.org/opportunity/global-talent/1329895", "Office Manager (EU only)", "Hamburg, Deutschland", "No", "66 applicants", "6 - 18 Months", "Caps & Collars GmbH"),
    @("1326661", "https://aiesec.org/opportunity/global-talent/1326661", "AL & ML Intern", "Manipal, Karnataka, India", "No", "18 applicants", "9 - 12 Weeks", "M.A.H.E."),
    @("1310446", "https://aiesec.org/opportunity/global-talent/1310446", "Education Coordinator", "Bursa, Türkiye", "No", "40 applicants", "9 - 12 Weeks", "Genç Kardelen Kindergarden")
)

$rowIndex = 2
foreach ($rowData in $data) {
    # Column A holds a numeric-looking opportunity id that must stay TEXT.
    # Temporarily mark the cell as Text so Excel doesn't coerce the digits
    # into a number, then strip the formatting back off so the cell keeps
    # the workbook's default (unstyled) look, matching plain scraped data.
    $idCell = $ws.Cells.Item($rowIndex, 1)
    $idCell.NumberFormat = "@"
    $idCell.Value = $rowData[0]
    $idCell.ClearFormats()

    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowData[$col - 1]
    }
    $rowIndex++
}
